# Apply underline/strikethrough formatting to specific bullet paragraphs
# (hook-from-mouse edit: mark items already covered / still-needed with
# strikethrough / underline so at-a-glance status is visible).

$d = $word.ActiveDocument

function Format-ParagraphByText($SearchText, $Kind) {
    $range = $d.Content
    $found = $range.Find.Execute($SearchText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found) {
        # Grab the whole paragraph (including the paragraph mark) so the
        # formatting lands on both the run(s) and the paragraph mark run
        # properties, matching Word's "select paragraph, apply format" flow.
        $para = $range.Paragraphs(1)
        $paraRange = $para.Range
        if ($Kind -eq "underline") {
            $paraRange.Font.Underline = 1
        } else {
            $paraRange.Font.StrikeThrough = 1
        }
    }
}

Format-ParagraphByText "Show the amount of boost the player has with a slider" "underline"
Format-ParagraphByText "Manage the different background tiles, spawning them as the player gets past certain boundaries to make endless level" "strike"
Format-ParagraphByText "Has Boost in it that player can pickup" "strike"
Format-ParagraphByText "When player picks it up, they gain boost" "underline"
